$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9921649098396301
$ws.Range("B1").Value = 0.9040123224258423
$ws.Range("C1").Value = 3.761794328689575
$ws.Range("D1").Value = 2.834225177764893
$ws.Range("E1").Value = 1.288111925125122
